$wb = $excel.ActiveWorkbook

# Rename sheets: remove spaces -> underscores
$wb.Worksheets.Item("RITM Sheet").Name = "RITM_Sheet"
$wb.Worksheets.Item("INC Sheet").Name = "INC_Sheet"

# Update window size/position
$excel.Width = 18160
$excel.Height = 12300
$excel.Left = 2320
$excel.Top = 500

# Update selection on RITM_Sheet
$ws1 = $wb.Worksheets.Item("RITM_Sheet")
$ws1.Activate()
$ws1.Range("C2").Select() | Out-Null

# Update selection on INC_Sheet
$ws2 = $wb.Worksheets.Item("INC_Sheet")
$ws2.Activate()
$ws2.Range("B1:B1048576").Select() | Out-Null

# Re-activate first sheet as tabSelected
$ws1.Activate()
